# Updates the lifeline-label shapes in the sequence diagram:
#   "LoanCard:UI"  -> ":LoanCard"      (split into two runs: ":" + "LoanCard")
#   "Nric:Model"   -> "nric:Nric"
#   "Phone:Model"  -> "phone:Phone"
#   "Email:Model"  -> "email:Email"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) {
        continue
    }

    $tr = $shp.TextFrame.TextRange
    $text = $tr.Text

    if ($text -eq "LoanCard:UI") {
        # Drop the trailing ":UI" suffix, leaving the base run as "LoanCard",
        # then prepend a fresh ":" run in front of it.
        $tr.Characters(9, 3).Text = ""
        [void]$tr.InsertBefore(":")
    }
    elseif ($text -eq "Nric:Model") {
        $tr.Text = "nric:Nric"
    }
    elseif ($text -eq "Phone:Model") {
        $tr.Text = "phone:Phone"
    }
    elseif ($text -eq "Email:Model") {
        $tr.Text = "email:Email"
    }
}
